# "Trocar as secoes de mapa e analise de ocorrencias" no dashboard.
# Reposiciona/redimensiona os cartoes (roundRect) do layout: os cartoes
# da esquerda passam a comecar em x=123750 EMU (em vez de ~198000/207001)
# e os cartoes/barras da direita se estendem para preencher o espaco
# liberado, mantendo a borda direita onde estava (ou esticando ate a
# nova borda do slide).
#
# PowerPoint's Shape.Left/.Width (etc.) are expressed in points
# (1 pt = 12700 EMU) and are stored internally as single-precision
# floats, so a naive "emu / 12700" assignment can round-trip to one
# EMU below the intended target for some values. The literal point
# values below were chosen so that, after the float32 round-trip, the
# saved OOXML has exactly the target EMU values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Cartao superior esquerdo ("Retangulo: Cantos Arredondados 18", id=19)
# off: 198000,767131 ext: 2890498,1117950 -> off: 123750,767131 ext: 2964748,1117950
$sh = $s.Shapes.Item("Retângulo: Cantos Arredondados 18")
$sh.Left = 9.744094488188976
$sh.Width = 233.44473266601562

# Cartao do meio esquerdo ("Retangulo: Cantos Arredondados 22", id=23)
# off: 207001,2083081 ext: 5879998,2238907 -> off: 123750,2083081 ext: 5963249,2238907
$sh = $s.Shapes.Item("Retângulo: Cantos Arredondados 22")
$sh.Left = 9.744094488188976
$sh.Width = 469.54716535433073

# Cartao inferior esquerdo ("Retangulo: Cantos Arredondados 23", id=24)
# off: 198000,4421093 ext: 5879998,2238907 -> off: 123750,4421093 ext: 5954248,2238907
$sh = $s.Shapes.Item("Retângulo: Cantos Arredondados 23")
$sh.Left = 9.744094488188976
$sh.Width = 468.8384251968504

# Cartao do mapa, a direita ("Retangulo: Cantos Arredondados 24", id=25)
# off inalterado: 6204002,2057300 ext: 5789997,4622294 -> ext: 5864248,4622294
$sh = $s.Shapes.Item("Retângulo: Cantos Arredondados 24")
$sh.Width = 461.7518310546875

# Barra superior ("Retangulo: Cantos Arredondados 4", id=5)
# off inalterado: 0,0 ext: 12068250,588725 -> ext: 12192000,588725
$sh = $s.Shapes.Item("Retângulo: Cantos Arredondados 4")
$sh.Width = 960.0
